$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# optimization_parameters sheet: restructure the parameter table.
#   - drop the stray repeated "value" cells in C1:F1
#   - rename "Model" -> "production_function"
#   - insert a new "L_curve" parameter row right after it
#   - drop the old "Deletion" row (now just before simulation_timepoints)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 header row only spans A:B now.
$ws.Range("C1:F1").Clear()

# Insert a new row at position 9 for the "L_curve" parameter; this shifts the
# existing estimate_params..simulation_timepoints rows down by one.
$ws.Rows.Item(9).Insert()

# Row 8: "Model" -> "production_function".
$ws.Range("A8").Value = "production_function"

# Row 9 (new): "L_curve" parameter = 1.
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# Remove the old "Deletion" row (now shifted down to row 17, right before
# simulation_timepoints).
$ws.Rows.Item(17).Select()
$ws.Rows.Item(17).Delete()

Write-Host "optimization_parameters updated"
